$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.07"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.75"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.209"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06098"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.515"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.720"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.355"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7991"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1577"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08102"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09280"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.897"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001704"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04822"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006157"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006186"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001102"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003393"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.694"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.259"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3360"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006162"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04590"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1122"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003128"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003402"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01021"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006023"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7497"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05877"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01010"
